$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "2025-06-25 17:58:25"
$ws.Range("B66").Value = "Policy Iteration"
$ws.Range("C66").Value = "LineWorld"
$ws.Range("D66").Value = 1
$ws.Range("E66").Value = 1
$ws.Range("F66").Value = 0.99
$ws.Cells.Item(66, 7).Value = "'"
$ws.Cells.Item(66, 7).Style = "Normal"
$ws.Cells.Item(66, 8).Value = "'"
$ws.Cells.Item(66, 8).Style = "Normal"
$ws.Cells.Item(66, 9).Value = "'"
$ws.Cells.Item(66, 9).Style = "Normal"

$ws.Range("A67").Value = "2025-06-25 17:58:29"
$ws.Range("B67").Value = "Policy Iteration"
$ws.Range("C67").Value = "LineWorld"
$ws.Range("D67").Value = 1
$ws.Range("E67").Value = 2
$ws.Range("F67").Value = 0.99
$ws.Cells.Item(67, 7).Value = "'"
$ws.Cells.Item(67, 7).Style = "Normal"
$ws.Cells.Item(67, 8).Value = "'"
$ws.Cells.Item(67, 8).Style = "Normal"
$ws.Cells.Item(67, 9).Value = "'"
$ws.Cells.Item(67, 9).Style = "Normal"

$ws.Range("A68").Value = "2025-06-25 17:58:31"
$ws.Range("B68").Value = "Policy Iteration"
$ws.Range("C68").Value = "LineWorld"
$ws.Range("D68").Value = 1
$ws.Range("E68").Value = 3
$ws.Range("F68").Value = 0.99
$ws.Cells.Item(68, 7).Value = "'"
$ws.Cells.Item(68, 7).Style = "Normal"
$ws.Cells.Item(68, 8).Value = "'"
$ws.Cells.Item(68, 8).Style = "Normal"
$ws.Cells.Item(68, 9).Value = "'"
$ws.Cells.Item(68, 9).Style = "Normal"

$ws.Range("A69").Value = "2025-06-25 17:59:29"
$ws.Range("B69").Value = "Policy Iteration"
$ws.Range("C69").Value = "GridWorld"
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = 1
$ws.Range("F69").Value = 0.99
$ws.Cells.Item(69, 7).Value = "'"
$ws.Cells.Item(69, 7).Style = "Normal"
$ws.Cells.Item(69, 8).Value = "'"
$ws.Cells.Item(69, 8).Style = "Normal"
$ws.Cells.Item(69, 9).Value = "'"
$ws.Cells.Item(69, 9).Style = "Normal"

$ws.Range("A70").Value = "2025-06-25 17:59:36"
$ws.Range("B70").Value = "Policy Iteration"
$ws.Range("C70").Value = "GridWorld"
$ws.Range("D70").Value = 1
$ws.Range("E70").Value = 2
$ws.Range("F70").Value = 0.99
$ws.Cells.Item(70, 7).Value = "'"
$ws.Cells.Item(70, 7).Style = "Normal"
$ws.Cells.Item(70, 8).Value = "'"
$ws.Cells.Item(70, 8).Style = "Normal"
$ws.Cells.Item(70, 9).Value = "'"
$ws.Cells.Item(70, 9).Style = "Normal"

$ws.Range("A71").Value = "2025-06-25 17:59:39"
$ws.Range("B71").Value = "Policy Iteration"
$ws.Range("C71").Value = "GridWorld"
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 3
$ws.Range("F71").Value = 0.99
$ws.Cells.Item(71, 7).Value = "'"
$ws.Cells.Item(71, 7).Style = "Normal"
$ws.Cells.Item(71, 8).Value = "'"
$ws.Cells.Item(71, 8).Style = "Normal"
$ws.Cells.Item(71, 9).Value = "'"
$ws.Cells.Item(71, 9).Style = "Normal"

$ws.Range("A72").Value = "2025-06-25 19:16:12"
$ws.Range("B72").Value = "Policy Iteration"
$ws.Range("C72").Value = "MontyHall LV1"
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 1
$ws.Range("F72").Value = 0.99
$ws.Cells.Item(72, 7).Value = "'"
$ws.Cells.Item(72, 7).Style = "Normal"
$ws.Cells.Item(72, 8).Value = "'"
$ws.Cells.Item(72, 8).Style = "Normal"
$ws.Cells.Item(72, 9).Value = "'"
$ws.Cells.Item(72, 9).Style = "Normal"

$ws.Range("A73").Value = "2025-06-25 19:16:16"
$ws.Range("B73").Value = "Policy Iteration"
$ws.Range("C73").Value = "MontyHall LV1"
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 2
$ws.Range("F73").Value = 0.99
$ws.Cells.Item(73, 7).Value = "'"
$ws.Cells.Item(73, 7).Style = "Normal"
$ws.Cells.Item(73, 8).Value = "'"
$ws.Cells.Item(73, 8).Style = "Normal"
$ws.Cells.Item(73, 9).Value = "'"
$ws.Cells.Item(73, 9).Style = "Normal"

$ws.Range("A74").Value = "2025-06-25 19:16:20"
$ws.Range("B74").Value = "Policy Iteration"
$ws.Range("C74").Value = "MontyHall LV1"
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 3
$ws.Range("F74").Value = 0.99
$ws.Cells.Item(74, 7).Value = "'"
$ws.Cells.Item(74, 7).Style = "Normal"
$ws.Cells.Item(74, 8).Value = "'"
$ws.Cells.Item(74, 8).Style = "Normal"
$ws.Cells.Item(74, 9).Value = "'"
$ws.Cells.Item(74, 9).Style = "Normal"

$ws.Range("A75").Value = "2025-06-25 19:36:53"
$ws.Range("B75").Value = "Value Iteration"
$ws.Range("C75").Value = "LineWorld"
$ws.Range("D75").Value = 1
$ws.Range("E75").Value = 1
$ws.Range("F75").Value = 0.99
$ws.Cells.Item(75, 7).Value = "'"
$ws.Cells.Item(75, 7).Style = "Normal"
$ws.Cells.Item(75, 8).Value = "'"
$ws.Cells.Item(75, 8).Style = "Normal"
$ws.Cells.Item(75, 9).Value = "'"
$ws.Cells.Item(75, 9).Style = "Normal"

$ws.Range("A76").Value = "2025-06-25 19:36:57"
$ws.Range("B76").Value = "Value Iteration"
$ws.Range("C76").Value = "LineWorld"
$ws.Range("D76").Value = 1
$ws.Range("E76").Value = 2
$ws.Range("F76").Value = 0.99
$ws.Cells.Item(76, 7).Value = "'"
$ws.Cells.Item(76, 7).Style = "Normal"
$ws.Cells.Item(76, 8).Value = "'"
$ws.Cells.Item(76, 8).Style = "Normal"
$ws.Cells.Item(76, 9).Value = "'"
$ws.Cells.Item(76, 9).Style = "Normal"

$ws.Range("A77").Value = "2025-06-25 19:36:59"
$ws.Range("B77").Value = "Value Iteration"
$ws.Range("C77").Value = "LineWorld"
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 3
$ws.Range("F77").Value = 0.99
$ws.Cells.Item(77, 7).Value = "'"
$ws.Cells.Item(77, 7).Style = "Normal"
$ws.Cells.Item(77, 8).Value = "'"
$ws.Cells.Item(77, 8).Style = "Normal"
$ws.Cells.Item(77, 9).Value = "'"
$ws.Cells.Item(77, 9).Style = "Normal"

$ws.Range("A78").Value = "2025-06-25 19:38:47"
$ws.Range("B78").Value = "Value Iteration"
$ws.Range("C78").Value = "GridWorld"
$ws.Range("D78").Value = 1
$ws.Range("E78").Value = 1
$ws.Range("F78").Value = 0.99
$ws.Cells.Item(78, 7).Value = "'"
$ws.Cells.Item(78, 7).Style = "Normal"
$ws.Cells.Item(78, 8).Value = "'"
$ws.Cells.Item(78, 8).Style = "Normal"
$ws.Cells.Item(78, 9).Value = "'"
$ws.Cells.Item(78, 9).Style = "Normal"

$ws.Range("A79").Value = "2025-06-25 19:38:53"
$ws.Range("B79").Value = "Value Iteration"
$ws.Range("C79").Value = "GridWorld"
$ws.Range("D79").Value = 1
$ws.Range("E79").Value = 2
$ws.Range("F79").Value = 0.99
$ws.Cells.Item(79, 7).Value = "'"
$ws.Cells.Item(79, 7).Style = "Normal"
$ws.Cells.Item(79, 8).Value = "'"
$ws.Cells.Item(79, 8).Style = "Normal"
$ws.Cells.Item(79, 9).Value = "'"
$ws.Cells.Item(79, 9).Style = "Normal"

$ws.Range("A80").Value = "2025-06-25 19:39:03"
$ws.Range("B80").Value = "Value Iteration"
$ws.Range("C80").Value = "GridWorld"
$ws.Range("D80").Value = 1
$ws.Range("E80").Value = 3
$ws.Range("F80").Value = 0.99
$ws.Cells.Item(80, 7).Value = "'"
$ws.Cells.Item(80, 7).Style = "Normal"
$ws.Cells.Item(80, 8).Value = "'"
$ws.Cells.Item(80, 8).Style = "Normal"
$ws.Cells.Item(80, 9).Value = "'"
$ws.Cells.Item(80, 9).Style = "Normal"

$ws.Range("A81").Value = "2025-06-25 19:39:09"
$ws.Range("B81").Value = "Value Iteration"
$ws.Range("C81").Value = "GridWorld"
$ws.Range("D81").Value = 1
$ws.Range("E81").Value = 4
$ws.Range("F81").Value = 0.99
$ws.Cells.Item(81, 7).Value = "'"
$ws.Cells.Item(81, 7).Style = "Normal"
$ws.Cells.Item(81, 8).Value = "'"
$ws.Cells.Item(81, 8).Style = "Normal"
$ws.Cells.Item(81, 9).Value = "'"
$ws.Cells.Item(81, 9).Style = "Normal"

$ws.Range("A82").Value = "2025-06-25 19:41:35"
$ws.Range("B82").Value = "Value Iteration"
$ws.Range("C82").Value = "MontyHall LV1"
$ws.Range("D82").Value = 0
$ws.Range("E82").Value = 1
$ws.Range("F82").Value = 0.99
$ws.Cells.Item(82, 7).Value = "'"
$ws.Cells.Item(82, 7).Style = "Normal"
$ws.Cells.Item(82, 8).Value = "'"
$ws.Cells.Item(82, 8).Style = "Normal"
$ws.Cells.Item(82, 9).Value = "'"
$ws.Cells.Item(82, 9).Style = "Normal"

$ws.Range("A83").Value = "2025-06-25 19:41:39"
$ws.Range("B83").Value = "Value Iteration"
$ws.Range("C83").Value = "MontyHall LV1"
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 2
$ws.Range("F83").Value = 0.99
$ws.Cells.Item(83, 7).Value = "'"
$ws.Cells.Item(83, 7).Style = "Normal"
$ws.Cells.Item(83, 8).Value = "'"
$ws.Cells.Item(83, 8).Style = "Normal"
$ws.Cells.Item(83, 9).Value = "'"
$ws.Cells.Item(83, 9).Style = "Normal"

$ws.Range("A84").Value = "2025-06-25 19:41:43"
$ws.Range("B84").Value = "Value Iteration"
$ws.Range("C84").Value = "MontyHall LV1"
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 3
$ws.Range("F84").Value = 0.99
$ws.Cells.Item(84, 7).Value = "'"
$ws.Cells.Item(84, 7).Style = "Normal"
$ws.Cells.Item(84, 8).Value = "'"
$ws.Cells.Item(84, 8).Style = "Normal"
$ws.Cells.Item(84, 9).Value = "'"
$ws.Cells.Item(84, 9).Style = "Normal"

$ws.Range("A85").Value = "2025-06-25 19:41:46"
$ws.Range("B85").Value = "Value Iteration"
$ws.Range("C85").Value = "MontyHall LV1"
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = 0.99
$ws.Cells.Item(85, 7).Value = "'"
$ws.Cells.Item(85, 7).Style = "Normal"
$ws.Cells.Item(85, 8).Value = "'"
$ws.Cells.Item(85, 8).Style = "Normal"
$ws.Cells.Item(85, 9).Value = "'"
$ws.Cells.Item(85, 9).Style = "Normal"

$ws.Range("A86").Value = "2025-06-25 19:41:49"
$ws.Range("B86").Value = "Value Iteration"
$ws.Range("C86").Value = "MontyHall LV1"
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 5
$ws.Range("F86").Value = 0.99
$ws.Cells.Item(86, 7).Value = "'"
$ws.Cells.Item(86, 7).Style = "Normal"
$ws.Cells.Item(86, 8).Value = "'"
$ws.Cells.Item(86, 8).Style = "Normal"
$ws.Cells.Item(86, 9).Value = "'"
$ws.Cells.Item(86, 9).Style = "Normal"

